# The workbook has a "Sheet1" worksheet with a helper formula in column A
# (rows 1-193) that builds a docker/bombardier command line from the IP
# address stored in column B, e.g.:
#   ="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B95&"&& sleep 5;"
#
# Rows 95 through 193 have an empty B (and C) column, so their A-column
# formula evaluates to a "template" string with a blank IP address. The
# commit removes that now-unused helper formula (and its cached/ shared
# string result) from rows 95-193, leaving those A cells blank while
# keeping their existing cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the formula/value from A95:A193, preserving cell formatting (style).
$ws.Range("A95:A193").ClearContents()
